$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "worker 1"
$ws.Range("B2").Value = "Mon 13-14, Tue 13-14, Wed 13-14, Wed 14-15, Wed 15-16, Thu 13-14, Thu 14-15, Fri 12-13, Fri 13-14, Fri 14-15, Fri 15-16, Fri 16-17"

$ws.Range("A3").Value = "worker 2"

$ws.Range("A4").Value = "worker 3"
$ws.Range("B4").Value = "Mon 14-15, Tue 12-13, Tue 16-17, Wed 8-9, Wed 9-10, Wed 12-13, Wed 14-15, Wed 15-16, Thu 12-13, Thu 16-17, Fri 9-10, Fri 12-13"

$ws.Range("A5").Value = "worker 4"
$ws.Range("B5").Value = "Mon 9-10, Mon 12-13, Mon 14-15, Mon 15-16, Mon 16-17, Tue 14-15, Wed 12-13, Wed 13-14, Thu 14-15, Thu 15-16, Fri 12-13, Fri 15-16"

$ws.Range("A6").Value = "worker 5"
$ws.Range("B6").Value = "Tue 15-16, Fri 11-12, Fri 13-14"

$ws.Range("A7").Value = "worker 6"
$ws.Range("B7").Value = "Tue 9-10, Tue 13-14, Thu 9-10, Thu 13-14, Fri 9-10, Fri 10-11, Fri 12-13, Fri 13-14"

$ws.Range("A8").Value = "worker 7"
$ws.Range("B8").Value = "Mon 12-13, Mon 13-14, Mon 15-16, Tue 12-13, Tue 13-14, Wed 12-13, Wed 13-14, Wed 14-15"

$ws.Range("A9").Value = "worker 8"
$ws.Range("B9").Value = "Mon 10-11, Tue 12-13, Tue 15-16, Wed 10-11, Wed 11-12, Wed 15-16, Thu 12-13, Thu 16-17, Fri 15-16, Fri 16-17"

$ws.Range("A10").Value = "worker 9"
$ws.Range("B10").Value = "Mon 9-10, Mon 10-11, Tue 14-15, Tue 15-16, Wed 9-10, Wed 10-11, Wed 16-17, Thu 14-15, Thu 15-16, Fri 10-11, Fri 12-13, Fri 13-14"

$ws.Range("A11").Value = "worker 10"
$ws.Range("B11").Value = "Mon 13-14, Wed 13-14, Thu 13-14, Fri 13-14, Fri 14-15"

$ws.Range("A12").Value = "worker 11"
$ws.Range("B12").Value = "Mon 12-13, Mon 13-14, Mon 14-15, Tue 12-13, Tue 13-14, Tue 15-16, Tue 16-17, Wed 13-14, Wed 16-17, Thu 12-13, Thu 13-14, Thu 14-15, Thu 15-16, Fri 11-12, Fri 14-15"

$ws.Range("A13").Value = "worker 12"
$ws.Range("B13").Value = "Mon 15-16, Tue 10-11, Tue 14-15, Wed 14-15, Fri 14-15"

$ws.Range("A14").Value = "worker 13"
$ws.Range("B14").Value = "Mon 8-9, Mon 11-12, Tue 8-9, Tue 13-14, Tue 14-15, Wed 8-9, Thu 8-9, Fri 8-9"

$ws.Range("A15").Value = "worker 14"
$ws.Range("B15").Value = "Mon 11-12, Mon 12-13, Mon 16-17, Tue 16-17, Wed 12-13, Thu 12-13, Thu 13-14, Thu 15-16, Thu 16-17"

$ws.Range("A16").Value = "worker 15"
$ws.Range("B16").Value = "Mon 8-9, Tue 8-9, Wed 8-9, Wed 11-12, Wed 12-13, Thu 8-9, Fri 8-9"

$ws.Range("A17").Value = "worker 16"
$ws.Range("B17").Value = "Mon 12-13, Mon 13-14, Mon 14-15, Thu 10-11"
